# Refresh the "cryptos" price list: update Price (D) and Volume(1h) (E)
# columns for the latest snapshot. Many Price values look like plain
# decimals ("8.43", "108.93", ...) which Excel would silently coerce to
# numbers on a normal `.Value =` assignment (losing the source formatting
# / turning them into binary floats like 8.4299999999999997). To keep
# them as literal text - matching how the sheet already stores every
# other cell in these columns (inline/shared strings) - we temporarily
# force the cell to Text format before writing the value, then clear the
# number format again so the cell's style stays at its original default.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "56.944.42"
$ws.Range("E2").Value = "  +4.35%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.264.38"
$ws.Range("E3").Value = "  +2.78%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "396.47"
$ws.Range("E5").Value = "  -1.21%  "

# Row 6 - Solana
Set-TextValue "D6" "108.93"
$ws.Range("E6").Value = "  -0.53%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +6.61%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.262.94"
$ws.Range("E8").Value = "  +2.86%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.01%  "

# Row 10 - Cardano
Set-TextValue "D10" "0.629"
$ws.Range("E10").Value = "  +1.90%  "

# Row 11 - Avalanche
$ws.Range("E11").Value = "  +0.50%  "

# Row 12 - Dogecoin
Set-TextValue "D12" "0.0980"
$ws.Range("E12").Value = "  +10.47%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +2.18%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.778.40"
$ws.Range("E14").Value = "  +2.69%  "

# Row 15 - Polkadot
Set-TextValue "D15" "8.43"
$ws.Range("E15").Value = "  +4.27%  "

# Row 16 - Chainlink
$ws.Range("E16").Value = "  +0.66%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.267.37"
$ws.Range("E17").Value = "  +2.80%  "

# Row 18 - Polygon
$ws.Range("E18").Value = "  -2.62%  "

# Row 19 - Uniswap
Set-TextValue "D19" "10.78"
$ws.Range("E19").Value = "  +2.29%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "56.852.01"
$ws.Range("E20").Value = "  +4.28%  "

# Row 21 - ImmutableX
$ws.Range("E21").Value = "  +1.72%  "

# Row 22 - ShibaInu
$ws.Range("E22").Value = "  +8.10%  "

# Row 23 - InternetComputer(DFINITY)
Set-TextValue "D23" "13.06"
$ws.Range("E23").Value = "  +1.29%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "295.86"
$ws.Range("E24").Value = "  +7.49%  "

# Row 25 - Litecoin
Set-TextValue "D25" "74.34"
$ws.Range("E25").Value = "  +2.85%  "

# Row 26 - PancakeSwap
$ws.Range("E26").Value = "  -2.49%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  +1.20%  "

# Row 28 - LEO
Set-TextValue "D28" "4.36"
$ws.Range("E28").Value = "  +1.00%  "

# Row 29 - Filecoin
$ws.Range("E29").Value = "  -3.68%  "

# Row 30 - RenderToken
Set-TextValue "D30" "7.27"
$ws.Range("E30").Value = "  -4.31%  "

# Row 31 - Kaspa
$ws.Range("E31").Value = "  +0.08%  "

# Row 32 - Dai
$ws.Range("E32").Value = "  +0.07%  "

# Row 33 - Cosmos
Set-TextValue "D33" "11.27"
$ws.Range("E33").Value = "  +2.20%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  -2.95%  "

# Row 35 - InjectiveProtocol
Set-TextValue "D35" "40.13"
$ws.Range("E35").Value = "  +9.18%  "

# Row 36 - VeChain
Set-TextValue "D36" "0.0493"
$ws.Range("E36").Value = "  -3.31%  "

# Row 37 - Toncoin
$ws.Range("E37").Value = "  +0.98%  "

# Row 38 - OKB
Set-TextValue "D38" "51.51"
$ws.Range("E38").Value = "  +0.05%  "

# Row 39 - FirstDigitalUSD
$ws.Range("E39").Value = "  -0.10%  "

# Row 40 - LidoDAOToken
$ws.Range("E40").Value = "  -3.75%  "

# Row 41 - Stacks
Set-TextValue "D41" "2.95"
$ws.Range("E41").Value = "  +2.03%  "

# Row 42 - Monero
Set-TextValue "D42" "139.05"
$ws.Range("E42").Value = "  +5.89%  "

# Row 43 - Stellar
$ws.Range("E43").Value = "  +3.89%  "

# Row 44 - NEARProtocol
$ws.Range("E44").Value = "  -1.72%  "

# Row 45 - ARBITRUM
$ws.Range("E45").Value = "  -1.88%  "

# Row 46 - Celestia
Set-TextValue "D46" "17.14"
$ws.Range("E46").Value = "  -0.78%  "

# Row 47 - TheGraph
$ws.Range("E47").Value = "  -3.15%  "

# Row 48 - EnergySwap
Set-TextValue "D48" "22.26"
$ws.Range("E48").Value = "  +0.71%  "

# Row 49 - WEMIXToken
$ws.Range("E49").Value = "  +3.39%  "

# Row 50 - Maker
$ws.Range("D50").Value = "2.163.97"
$ws.Range("E50").Value = "  +3.42%  "

# Row 51 - ApeXProtocol
$ws.Range("E51").Value = "  -5.23%  "
